# "What we did in class 10/7": re-sort the Learning Outcomes table (Table1)
# by Area then Course (ascending), matching Data > Sort on the table, and
# leave the selection where it ended up afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sortRange = $ws.Range("A1:D142")
$key1 = $ws.Range("A2:A142")
$key2 = $ws.Range("B2:B142")

# xlAscending = 1, xlYes (has headers) = 1
$sortRange.Sort($key1, 1, $key2, $null, 1, $null, $null, 1) | Out-Null

$ws.Range("C128").Select() | Out-Null
